# Applies updated market-price/profit figures to the crafting-class sheets
# (values recomputed by the scheduled market data runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H9").Value = 294.60715
$ws.Range("J9").Value = 351
$ws.Range("L9").Value = 351
$ws.Range("N9").Value = -689

$ws.Range("H17").Value = 1680211.9
$ws.Range("I17").Value = 900
$ws.Range("J17").Value = 1832876.5
$ws.Range("K17").Value = 2700
$ws.Range("L17").Value = 5498629.5
$ws.Range("M17").Value = -2532
$ws.Range("N17").Value = -5498965.5

$ws.Range("H40").Value = 96914.47
$ws.Range("I40").Value = 601076.2
$ws.Range("J40").Value = 3551.1853
$ws.Range("K40").Value = 601076.2
$ws.Range("L40").Value = 3551.1853
$ws.Range("M40").Value = -600901.2
$ws.Range("N40").Value = -3901.1853

$ws.Range("H62").Value = 9648.5
$ws.Range("J62").Value = 10497
$ws.Range("L62").Value = 10497
$ws.Range("N62").Value = -11745

$ws.Range("H65").Value = 9648.5
$ws.Range("J65").Value = 10497
$ws.Range("L65").Value = 52485
$ws.Range("N65").Value = -58725

$ws.Range("H70").Value = 50002076
$ws.Range("I70").Value = 4998
$ws.Range("J70").Value = 62501348
$ws.Range("K70").Value = 14994
$ws.Range("L70").Value = 187504044
$ws.Range("M70").Value = -14724
$ws.Range("N70").Value = -187504584

$ws.Range("H73").Value = 50002076
$ws.Range("I73").Value = 4998
$ws.Range("J73").Value = 62501348
$ws.Range("K73").Value = 14994
$ws.Range("L73").Value = 187504044
$ws.Range("M73").Value = -14058
$ws.Range("N73").Value = -187505916

$ws.Range("H98").Value = 922.65717
$ws.Range("I98").Value = 854.931
$ws.Range("J98").Value = 1250
$ws.Range("K98").Value = 854.931
$ws.Range("L98").Value = 1250
$ws.Range("M98").Value = 643.069
$ws.Range("N98").Value = -4246

$ws.Range("H99").Value = 83334180
$ws.Range("I99").Value = 111112170
$ws.Range("K99").Value = 333336510
$ws.Range("M99").Value = -333335012

$ws.Range("H122").Value = 922.65717
$ws.Range("I122").Value = 854.931
$ws.Range("J122").Value = 1250
$ws.Range("K122").Value = 2564.793
$ws.Range("L122").Value = 3750
$ws.Range("M122").Value = -114.7930000000001
$ws.Range("N122").Value = -8650

$ws.Range("H125").Value = 2013.0286
$ws.Range("I125").Value = 1163.4615
$ws.Range("J125").Value = 2515.0454
$ws.Range("K125").Value = 10471.1535
$ws.Range("L125").Value = 22635.4086
$ws.Range("M125").Value = -8011.153499999999
$ws.Range("N125").Value = -27555.4086

$ws.Range("H135").Value = 1230.4857
$ws.Range("I135").Value = 1180.5
$ws.Range("J135").Value = 2930
$ws.Range("K135").Value = 10624.5
$ws.Range("L135").Value = 26370
$ws.Range("M135").Value = -8089.5
$ws.Range("N135").Value = -31440

$ws.Range("H137").Value = 34215.31
$ws.Range("I137").Value = 48443.832
$ws.Range("J137").Value = 2201.125
$ws.Range("K137").Value = 145331.496
$ws.Range("L137").Value = 6603.375
$ws.Range("M137").Value = -142781.496
$ws.Range("N137").Value = -11703.375

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 20017558
$ws.Range("I32").Value = 21447992
$ws.Range("J32").Value = 10208866
$ws.Range("K32").Value = 21447992
$ws.Range("L32").Value = 10208866
$ws.Range("M32").Value = -21447705
$ws.Range("N32").Value = -10209440

$ws.Range("H61").Value = 3775.7334
$ws.Range("I61").Value = 3688.2856
$ws.Range("K61").Value = 3688.2856
$ws.Range("M61").Value = -3476.2856

$ws.Range("H74").Value = 1862
$ws.Range("I74").Value = 1915.75
$ws.Range("K74").Value = 1915.75
$ws.Range("M74").Value = -1041.75

$ws.Range("H77").Value = 1862
$ws.Range("I77").Value = 1915.75
$ws.Range("K77").Value = 9578.75
$ws.Range("M77").Value = -5210.75

$ws.Range("H102").Value = 1425.5
$ws.Range("I102").Value = 1260.65
$ws.Range("J102").Value = 2249.75
$ws.Range("K102").Value = 1260.65
$ws.Range("L102").Value = 2249.75
$ws.Range("M102").Value = 361.3499999999999
$ws.Range("N102").Value = -5493.75

$ws.Range("H110").Value = 2160.75
$ws.Range("J110").Value = 2650
$ws.Range("L110").Value = 2650
$ws.Range("N110").Value = -6740

$ws.Range("H122").Value = 5041.222
$ws.Range("I122").Value = 3437.0715
$ws.Range("K122").Value = 10311.2145
$ws.Range("M122").Value = -7861.2145

$ws.Range("H136").Value = 3775.7334
$ws.Range("I136").Value = 3688.2856
$ws.Range("K136").Value = 11064.8568
$ws.Range("M136").Value = -8514.856800000001

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H7").Value = 62.416668
$ws.Range("I7").Value = 69.71429000000001
$ws.Range("J7").Value = 52.2
$ws.Range("K7").Value = 69.71429000000001
$ws.Range("L7").Value = 52.2
$ws.Range("M7").Value = 43.28570999999999
$ws.Range("N7").Value = -278.2

$ws.Range("H68").Value = 50493
$ws.Range("J68").Value = 68986
$ws.Range("L68").Value = 68986
$ws.Range("N68").Value = -70484

$ws.Range("H71").Value = 50493
$ws.Range("J71").Value = 68986
$ws.Range("L71").Value = 206958
$ws.Range("N71").Value = -214446

$ws.Range("H103").Value = 50921.46
$ws.Range("I103").Value = 8864
$ws.Range("K103").Value = 8864
$ws.Range("M103").Value = -7692

$ws.Range("H105").Value = 1410.5834
$ws.Range("I105").Value = 1057.8422
$ws.Range("K105").Value = 1057.8422
$ws.Range("M105").Value = 689.1578

$ws.Range("H132").Value = 73540.73
$ws.Range("I132").Value = 46215.332
$ws.Range("J132").Value = 204702.6
$ws.Range("K132").Value = 138645.996
$ws.Range("L132").Value = 614107.8
$ws.Range("M132").Value = -136115.996
$ws.Range("N132").Value = -619167.8

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("N102").Value = 0
$ws.Range("L102").ClearContents()

$ws.Range("H116").Value = 1058.25
$ws.Range("J116").Value = 1194.3334
$ws.Range("L116").Value = 3583.0002
$ws.Range("N116").Value = -10467.0002

$ws.Range("H129").Value = 1779.4706
$ws.Range("I129").Value = 829.8889
$ws.Range("K129").Value = 2489.6667
$ws.Range("M129").Value = 2510.3333

$ws.Range("H137").Value = 4823.25
$ws.Range("I137").Value = 5646.5
$ws.Range("J137").Value = 4000
$ws.Range("K137").Value = 16939.5
$ws.Range("L137").Value = 12000
$ws.Range("M137").Value = -11839.5
$ws.Range("N137").Value = -22200

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H80").Value = 2647
$ws.Range("J80").Value = 3600
$ws.Range("L80").Value = 3600
$ws.Range("N80").Value = -5596

$ws.Range("H83").Value = 2647
$ws.Range("J83").Value = 3600
$ws.Range("L83").Value = 18000
$ws.Range("N83").Value = -27984

$ws.Range("H97").Value = 976.6
$ws.Range("I97").Value = 867.4516
$ws.Range("J97").Value = 1352.5555
$ws.Range("K97").Value = 867.4516
$ws.Range("L97").Value = 1352.5555
$ws.Range("M97").Value = -371.4516
$ws.Range("N97").Value = -2344.5555

$ws.Range("H99").Value = 39170.75
$ws.Range("I99").Value = 29480.857
$ws.Range("K99").Value = 29480.857
$ws.Range("M99").Value = -27234.857

$ws.Range("H102").Value = 1799.8
$ws.Range("I102").Value = 1624.75
$ws.Range("K102").Value = 1624.75
$ws.Range("M102").Value = -2.75

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H22").Value = 2092.1738
$ws.Range("J22").Value = 2544
$ws.Range("L22").Value = 2544
$ws.Range("N22").Value = -3134

$ws.Range("H27").Value = 2092.1738
$ws.Range("J27").Value = 2544
$ws.Range("L27").Value = 2544
$ws.Range("N27").Value = -2758

$ws.Range("H82").Value = 10879
$ws.Range("I82").Value = 2061.3076
$ws.Range("J82").Value = 23615.666
$ws.Range("K82").Value = 2061.3076
$ws.Range("L82").Value = 23615.666
$ws.Range("M82").Value = -1700.3076
$ws.Range("N82").Value = -24337.666

$ws.Range("H85").Value = 10879
$ws.Range("I85").Value = 2061.3076
$ws.Range("J85").Value = 23615.666
$ws.Range("K85").Value = 2061.3076
$ws.Range("L85").Value = 23615.666
$ws.Range("M85").Value = -813.3076000000001
$ws.Range("N85").Value = -26111.666

$ws.Range("H128").Value = 108979
$ws.Range("J128").Value = 108979
$ws.Range("L128").Value = 108979
$ws.Range("N128").Value = -118939

$ws.Range("H132").Value = 342668
$ws.Range("I132").Value = 342668
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 1028004
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = -1025474
$ws.Range("M132").ClearContents()

$ws.Range("H136").Value = 3011.8125
$ws.Range("I136").Value = 2464.182
$ws.Range("J136").Value = 4216.6
$ws.Range("K136").Value = 7392.545999999999
$ws.Range("L136").Value = 12649.8
$ws.Range("M136").Value = -4842.545999999999
$ws.Range("N136").Value = -17749.8

$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H18").Value = 3400.7

$ws.Range("H122").Value = 83339944
$ws.Range("I122").Value = 250004750
$ws.Range("J122").Value = 7547.625
$ws.Range("K122").Value = 750014250
$ws.Range("L122").Value = 22642.875
$ws.Range("M122").Value = -750011800
$ws.Range("N122").Value = -27542.875

$ws.Range("H126").Value = 7491.636

$ws.Range("H132").Value = 3980.8096
$ws.Range("I132").Value = 3888.6667
$ws.Range("K132").Value = 11666.0001
$ws.Range("M132").Value = -9136.000100000001

$ws.Range("H136").Value = 30773.8
$ws.Range("I136").Value = 2142.3333
$ws.Range("J136").Value = 52247.4
$ws.Range("K136").Value = 6426.999899999999
$ws.Range("L136").Value = 156742.2
$ws.Range("M136").Value = -3876.999899999999
$ws.Range("N136").Value = -161842.2
